# 3V3 boost converter replacement
# Fill in Farnell part numbers/prices and supplier choice for the inductors
# (L1, L2/L3) used by the new boost converter, rows 21 and 22 of Table1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21 - L1 (MH2029-300Y)
$ws.Range("E21").Value = 1515662
$ws.Range("F21").Value = 0.092
$ws.Range("Q21").Value = "Farnell"

# Row 22 - L2,L3 (4.7uH,1.2A)
$ws.Range("E22").Value = 2118126
$ws.Range("F22").Value = 0.26
$ws.Range("Q22").Value = "Farnell"

$wb.Save()
